$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Common" sheet (index 2): insert a new row at row 70 for a new
#    "VSTAT License File" entry, pushing all subsequent rows down by one.
# ---------------------------------------------------------------------------
$common = $wb.Worksheets.Item(2)

# The last used row before the insert is 129 (dimension A1:B129).
$lastRow = 129

$common.Rows("70:70").Insert()
$common.Range("A70").Value = "VSTAT License File"

# Row Insert() shifts cell values/styles down automatically, but cell
# comments stay anchored to their original row index. Re-home every comment
# from the old row N to the new row N+1, walking bottom-up so we never
# clobber a comment before it has been read.
for ($r = $lastRow; $r -ge 70; $r--) {
    $srcCell = $common.Cells.Item($r, 1)
    $cmt = $srcCell.Comment
    if ($cmt -ne $null) {
        $txt = $cmt.Text()
        $cmt.Delete()
        $dstCell = $common.Cells.Item($r + 1, 1)
        $dstCell.AddComment($txt)
    }
}

# Add the new comment for the newly inserted row.
$common.Cells.Item(70, 1).AddComment("Optional License File for Elasticsearch [default: ]")

# ---------------------------------------------------------------------------
# 2) "Credentials" sheet (index 7): reword a batch of existing comments and
#    one cell value (NETCONF Manager sudo password).
# ---------------------------------------------------------------------------
$cred = $wb.Worksheets.Item(7)

$cred.Range("A57").Value = "NETCONF Manager VM password for running sudo commands, and will be used for the installation of NETCONF Manager."

$cred.Range("A10").Comment.Text("VSD Username will be used for logging into VSD command line. Used for both Install and Upgrade procedures. [default: root]")
$cred.Range("A11").Comment.Text("VSD password will be used for logging into the command line. Used for both Install and Upgrade procedures. [default: Alcateldc]")
$cred.Range("A12").Comment.Text("VSC Username will be used for logging into command line (should have admin privileges). Used for upgrade procedure only [default: ]")
$cred.Range("A13").Comment.Text("VSC password will be used for logging into the command line. Used for upgrade procedure only [default: ]")
$cred.Range("A15").Comment.Text("ElasticSearch (Stats) Username will be used for logging into command line. Used for both Install and Upgrade procedures. [default: ]")
$cred.Range("A16").Comment.Text("ElasticSearch (Stats) password will be used for logging into the command line. Used for both Install and Upgrade procedures. [default: ]")
$cred.Range("A17").Comment.Text("ElasticSearch (Stats) root password required for VSTAT Upgrade only [default: ]")
$cred.Range("A19").Comment.Text("This VSD Username(also known as csproot user). Used for both Install and Upgrade procedures. Must have csproot privileges. [default: csproot]")
$cred.Range("A20").Comment.Text("This VSD password(also known as csproot password) will be used for API authentication. Used for both Install and Upgrade procedures. Must have csproot privileges. [default: csproot]")
$cred.Range("A21").Comment.Text("This VSD Mysql password. Used for both Install and Upgrade procedures. [default: ]")
$cred.Range("A37").Comment.Text("Username for OpenStack. [default: ]")
$cred.Range("A40").Comment.Text("vCenter Username. [default: ]")
$cred.Range("A43").Comment.Text("Username for Compute node to install VRS. [default: root]")
$cred.Range("A44").Comment.Text("Password for Compute node, and will be used for installation of VRS [default: ]")
$cred.Range("A54").Comment.Text("NFS username to login into command line, and will be used for NFS configuration. Default user is root. [default: root]")
$cred.Range("A56").Comment.Text("Username for NETCONF Manager VM, and will be used for the installation of NETCONF Manager. Default user is root. [default: root]")
$cred.Range("A58").Comment.Text("Username for NETCONF Manager user, and will be used for the installation of NETCONF Manager. [default: netconfmgr]")
$cred.Range("A59").Comment.Text("Password for NETCONF manager user, and will be used for the installation of NETCONF Manager. [default: password]")
$cred.Range("A61").Comment.Text("Username for SMTP Server, and will be used for Email health report.")
$cred.Range("A62").Comment.Text("Password for SMTP Server, and will be used for Email health report.")
$cred.Range("A64").Comment.Text("Username for the monit mail server.")
$cred.Range("A67").Comment.Text("Username for NUH notification application, and will be used for installation of NUH.")
$cred.Range("A68").Comment.Text("Password for NUH notification application, and will be used for installation of NUH.")
$cred.Range("A69").Comment.Text("Username for NUH notification application, and will be used for installation of NUH.")
$cred.Range("A70").Comment.Text("Password for NUH notification application, and will be used for installation of NUH.")
